$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 84.75

$ws.Range("H17").Value = 910.3137
$ws.Range("J17").Value = 894.4878
$ws.Range("L17").Value = 2683.4634
$ws.Range("N17").Value = -3019.4634

$ws.Range("H32").Value = 2000
$ws.Range("I32").Value = 2000
$ws.Range("J32").Value = 2000
$ws.Range("K32").Value = 2000
$ws.Range("L32").Value = 2000
$ws.Range("M32").Value = -1674
$ws.Range("N32").Value = -2652

$ws.Range("H112").Value = 3155.6667
$ws.Range("I112").Value = 3416.25
$ws.Range("J112").Value = 2947.2
$ws.Range("K112").Value = 10248.75
$ws.Range("L112").Value = 8841.599999999999
$ws.Range("M112").Value = -9140.75
$ws.Range("N112").Value = -11057.6

$ws.Range("H116").Value = 33943.5
$ws.Range("I116").Value = 18416.666
$ws.Range("K116").Value = 18416.666
$ws.Range("M116").Value = -14974.666

$ws.Range("H135").Value = 4518.5835
$ws.Range("J135").Value = 13000
$ws.Range("L135").Value = 117000
$ws.Range("N135").Value = -122070

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 310896.3
$ws.Range("I32").Value = 317486.8
$ws.Range("K32").Value = 317486.8
$ws.Range("M32").Value = -317199.8

$ws.Range("H74").Value = 5586.741
$ws.Range("I74").Value = 3114.8936
$ws.Range("J74").Value = 16148.272
$ws.Range("K74").Value = 3114.8936
$ws.Range("L74").Value = 16148.272
$ws.Range("M74").Value = -2240.8936
$ws.Range("N74").Value = -17896.272

$ws.Range("H77").Value = 5586.741
$ws.Range("I77").Value = 3114.8936
$ws.Range("J77").Value = 16148.272
$ws.Range("K77").Value = 15574.468
$ws.Range("L77").Value = 80741.36
$ws.Range("M77").Value = -11206.468
$ws.Range("N77").Value = -89477.36

$ws.Range("H110").Value = 1456
$ws.Range("I110").Value = 1361.909
$ws.Range("J110").Value = 1603.8572
$ws.Range("K110").Value = 1361.909
$ws.Range("L110").Value = 1603.8572
$ws.Range("M110").Value = 683.0909999999999
$ws.Range("N110").Value = -5693.8572

$ws.Range("H132").Value = 5266.647
$ws.Range("I132").Value = 3407.348
$ws.Range("K132").Value = 10222.044
$ws.Range("M132").Value = -7692.044

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H11").Value = 75
$ws.Range("I11").Value = 0
$ws.Range("K11").Value = 0
$ws.Range("M11").ClearContents()

$ws.Range("H134").Value = 3361.6667
$ws.Range("I134").Value = 3151.5
$ws.Range("K134").Value = 9454.5
$ws.Range("M134").Value = -6919.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2476.8
$ws.Range("I31").Value = 2728.2727
$ws.Range("K31").Value = 2728.2727
$ws.Range("M31").Value = -2433.2727

$ws.Range("H34").Value = 2476.8
$ws.Range("I34").Value = 2728.2727
$ws.Range("K34").Value = 2728.2727
$ws.Range("M34").Value = -2526.2727

$ws.Range("H58").Value = 6477.1724
$ws.Range("I58").Value = 5656
$ws.Range("J58").Value = 7243.6
$ws.Range("K58").Value = 5656
$ws.Range("L58").Value = 7243.6
$ws.Range("M58").Value = -5453
$ws.Range("N58").Value = -7649.6

$ws.Range("H99").Value = 20627.818
$ws.Range("I99").Value = 24461.777
$ws.Range("K99").Value = 24461.777
$ws.Range("M99").Value = -22963.777

$ws.Range("H105").Value = 6606.3
$ws.Range("I105").Value = 6740.3887
$ws.Range("J105").Value = 5399.5
$ws.Range("K105").Value = 6740.3887
$ws.Range("L105").Value = 5399.5
$ws.Range("M105").Value = -4993.3887
$ws.Range("N105").Value = -8893.5

$ws.Range("H126").Value = 20627.818
$ws.Range("I126").Value = 24461.777
$ws.Range("K126").Value = 73385.33099999999
$ws.Range("M126").Value = -70915.33099999999

$ws.Range("H132").Value = 2345.641
$ws.Range("I132").Value = 2313.4856
$ws.Range("K132").Value = 6940.4568
$ws.Range("M132").Value = -4410.4568

$ws.Range("H134").Value = 2719.923
$ws.Range("I134").Value = 2410.5293
$ws.Range("J134").Value = 3304.3333
$ws.Range("K134").Value = 7231.5879
$ws.Range("L134").Value = 9912.999899999999
$ws.Range("M134").Value = -4696.5879
$ws.Range("N134").Value = -14982.9999

$ws.Range("H136").Value = 6477.1724
$ws.Range("I136").Value = 5656
$ws.Range("J136").Value = 7243.6
$ws.Range("K136").Value = 16968
$ws.Range("L136").Value = 21730.8
$ws.Range("M136").Value = -14418
$ws.Range("N136").Value = -26830.8

$ws.Range("H141").Value = 49666.668
$ws.Range("J141").Value = 49666.668
$ws.Range("L141").Value = 49666.668
$ws.Range("N141").Value = -60026.668

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 182.8
$ws.Range("I2").Value = 89
$ws.Range("J2").Value = 206.25
$ws.Range("K2").Value = 534
$ws.Range("L2").Value = 1237.5
$ws.Range("M2").Value = -421
$ws.Range("N2").Value = -1463.5

$ws.Range("H14").Value = 94.818184
$ws.Range("I14").Value = 94.818184
$ws.Range("K14").Value = 284.454552
$ws.Range("M14").Value = -111.454552

$ws.Range("H68").Value = 5546
$ws.Range("I68").Value = 4359.6
$ws.Range("J68").Value = 6139.2
$ws.Range("K68").Value = 13078.8
$ws.Range("L68").Value = 18417.6
$ws.Range("M68").Value = -12267.8
$ws.Range("N68").Value = -20039.6

$ws.Range("H71").Value = 5546
$ws.Range("I71").Value = 4359.6
$ws.Range("J71").Value = 6139.2
$ws.Range("K71").Value = 39236.4
$ws.Range("L71").Value = 55252.8
$ws.Range("M71").Value = -35180.4
$ws.Range("N71").Value = -63364.8

$ws.Range("H87").Value = 9208
$ws.Range("I87").Value = 7811.5
$ws.Range("J87").Value = 12001
$ws.Range("K87").Value = 23434.5
$ws.Range("L87").Value = 36003
$ws.Range("N87").Value = -38499
$ws.Range("M87").Value = -22186.5

$ws.Range("H90").Value = 9208
$ws.Range("I90").Value = 7811.5
$ws.Range("J90").Value = 12001
$ws.Range("K90").Value = 70303.5
$ws.Range("L90").Value = 108009
$ws.Range("N90").Value = -120489
$ws.Range("M90").Value = -64063.5

$ws.Range("H110").Value = 8898.625
$ws.Range("I110").Value = 7884.143
$ws.Range("J110").Value = 16000
$ws.Range("K110").Value = 23652.429
$ws.Range("L110").Value = 48000
$ws.Range("M110").Value = -19562.429
$ws.Range("N110").Value = -56180

$ws.Range("H112").Value = 7374.75
$ws.Range("I112").Value = 999.5
$ws.Range("J112").Value = 13750
$ws.Range("K112").Value = 2998.5
$ws.Range("L112").Value = 41250
$ws.Range("M112").Value = -1890.5
$ws.Range("N112").Value = -43466

$ws.Range("H120").Value = 10338.167
$ws.Range("I120").Value = 10338.167
$ws.Range("J120").Value = 0
$ws.Range("K120").Value = 31014.501
$ws.Range("L120").Value = 0
$ws.Range("M120").Value = -26176.501
$ws.Range("N120").ClearContents()

$ws.Range("H122").Value = 734288.9
$ws.Range("J122").Value = 1288.6316
$ws.Range("L122").Value = 11597.6844
$ws.Range("N122").Value = -16497.6844

$ws.Range("H134").Value = 5629.722
$ws.Range("I134").Value = 1012.2727
$ws.Range("K134").Value = 3036.8181
$ws.Range("M134").Value = 2033.1819

$ws.Range("H137").Value = 2757.6
$ws.Range("J137").Value = 2849.4285
$ws.Range("L137").Value = 8548.2855
$ws.Range("N137").Value = -18748.2855

$ws.Range("H138").Value = 36483.3
$ws.Range("I138").Value = 64669.4
$ws.Range("J138").Value = 8297.200000000001
$ws.Range("K138").Value = 194008.2
$ws.Range("L138").Value = 24891.6
$ws.Range("M138").Value = -188868.2
$ws.Range("N138").Value = -35171.60000000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 11884.577
$ws.Range("I70").Value = 14605.211
$ws.Range("K70").Value = 14605.211
$ws.Range("M70").Value = -14335.211

$ws.Range("H73").Value = 11884.577
$ws.Range("I73").Value = 14605.211
$ws.Range("K73").Value = 14605.211
$ws.Range("M73").Value = -13669.211

$ws.Range("H97").Value = 957.4545000000001
$ws.Range("I97").Value = 920.8823
$ws.Range("J97").Value = 1081.8
$ws.Range("K97").Value = 920.8823
$ws.Range("L97").Value = 1081.8
$ws.Range("M97").Value = -424.8823
$ws.Range("N97").Value = -2073.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1884.8572
$ws.Range("I16").Value = 1918.6
$ws.Range("K16").Value = 1918.6
$ws.Range("M16").Value = -1748.6

$ws.Range("H55").Value = 1459.6451
$ws.Range("I55").Value = 1101.1333
$ws.Range("K55").Value = 1101.1333
$ws.Range("M55").Value = -928.1333

$ws.Range("H93").Value = 2573
$ws.Range("I93").Value = 1167.8182
$ws.Range("K93").Value = 1167.8182
$ws.Range("M93").Value = 80.18180000000007

$ws.Range("H122").Value = 5108.1665
$ws.Range("I122").Value = 3851
$ws.Range("K122").Value = 11553
$ws.Range("M122").Value = -9103

$ws.Range("H136").Value = 2821.25
$ws.Range("I136").Value = 2676
$ws.Range("K136").Value = 8028
$ws.Range("M136").Value = -5478

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2261.0334
$ws.Range("I132").Value = 1197.8889
$ws.Range("K132").Value = 3593.6667
$ws.Range("M132").Value = -1063.6667

$ws.Range("H136").Value = 1542.8334
$ws.Range("I136").Value = 1487.2
$ws.Range("K136").Value = 4461.6
$ws.Range("M136").Value = -1911.6
